$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# MySQL for Excel add-in's hidden helper name, (re)written on save.
$defName = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", "=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&"" ""&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)")
$defName.Visible = $false

# Replace the existing participant rows (2-4) with the new data, and add a
# new 5th row, per the 3/24/2018 8:11AM update.
$data = @(
    @("Gamboa",    "Rey Christian",      "Lopez"),
    @("Francisco", "Christopher Jorge",  "Pineda"),
    @("Naguit",    "Lanz",               "Pundavela"),
    @("Calantuan", "Earle",              "LeBron")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Column B is now wider to fit the longest new name ("Christopher Jorge").
$ws.Columns.Item(2).ColumnWidth = 15.75

# Match the saved selection from the authored workbook.
$null = $ws.Range("C10").Select()
